$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.321.55'
$ws.Range("E2").Value = '  +1.32%  '

$ws.Range("D3").Value = '1.570.54'
$ws.Range("E3").Value = '  +0.48%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = '211.10'
$ws.Range("E5").Value = '  +1.75%  '

$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").Value = '0.249'
$ws.Range("E9").Value = '  +0.51%  '

$ws.Range("E10").Value = '  +0.14%  '

$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("E12").Value = '  +0.68%  '

$ws.Range("D13").Value = '1.609.16'
$ws.Range("E13").Value = '  +3.05%  '

$ws.Range("E14").Value = '  +0.85%  '

$ws.Range("D15").Value = '0.520'
$ws.Range("E15").Value = '  +0.09%  '

$ws.Range("D16").Value = '27.261.01'
$ws.Range("E16").Value = '  +1.15%  '

$ws.Range("D17").Value = '62.31'
$ws.Range("E17").Value = '  +0.29%  '

$ws.Range("D18").Value = '7.53'
$ws.Range("E18").Value = '  +2.36%  '

$ws.Range("D19").Value = '217.76'
$ws.Range("E19").Value = '  +0.27%  '

$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("E21").Value = '  +0.03%  '

$ws.Range("E22").Value = '  +1.27%  '

$ws.Range("D23").Value = '9.24'
$ws.Range("E23").Value = '  +0.17%  '

$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("D25").Value = '153.49'
$ws.Range("E25").Value = '  +0.52%  '

$ws.Range("E26").Value = '  +0.63%  '

$ws.Range("D27").Value = '15.09'
$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("E28").Value = '  +1.93%  '

$ws.Range("E29").Value = '  +0.14%  '

$ws.Range("E30").Value = '  +2.61%  '

$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("E32").Value = '  +0.46%  '

$ws.Range("E33").Value = '  +1.69%  '

$ws.Range("D34").Value = '1.445.43'
$ws.Range("E34").Value = '  +1.72%  '

$ws.Range("E35").Value = '  +4.37%  '

$ws.Range("D36").Value = '1.62'
$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("E37").Value = '  +0.35%  '

$ws.Range("E38").Value = '  +0.94%  '

$ws.Range("D39").Value = '0.536'
$ws.Range("E39").Value = '  +0.55%  '

$ws.Range("D40").Value = '5.89'
$ws.Range("E40").Value = '  +2.41%  '

$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("E42").Value = '  +0.05%  '

$ws.Range("E43").Value = '  +0.50%  '

$ws.Range("E44").Value = '  -0.33%  '

$ws.Range("D45").Value = '64.69'
$ws.Range("E45").Value = '  -0.25%  '

$ws.Range("D46").Value = '1.73'
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").Value = '1.705.96'
$ws.Range("E47").Value = '  +0.58%  '

$ws.Range("D48").Value = '86.31'
$ws.Range("E48").Value = '  -1.41%  '

$ws.Range("D49").Value = '0.0527'
$ws.Range("E49").Value = '  +1.09%  '

$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").Value = '  +0.10%  '
